$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customers")

# Update Howard's birthday value (E5): 21/05/2002 -> 12/05/1987
$ws.Range("E5").Value = "12/05/1987"

# Turn Billy's canDrinkAlcohol cell (C6) into a =TRUE() formula instead of
# a literal boolean, keeping its existing number format/style.
$ws.Range("C6").Formula = "=TRUE()"

# Move/restore the active selection from E7 to E5.
$ws.Range("E5").Select()
